# Generate Report for handback
# Updates the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) for the newly handed-back
# file (row 2) on both the "zh-cn" and "de-de" report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-25 03:38:37"
$wsZhCn.Range("G2").Value = "2016-01-25 03:39:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-25 03:38:49"
$wsDeDe.Range("G2").Value = "2016-01-25 03:39:44"
